$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10: "2021年" in A10 (same formatting as the other year cells,
# e.g. A9) and its numeric value 1012 in B10.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A10").Value = "2021年"
$ws.Range("B10").Value = 1012
